# Apply the Mon Mar 25 22:47:16 UTC 2024 "Updated cryptos list" GitHub Actions refresh:
# latest prices / 1h volumes for every coin row, plus a few rank swaps where two
# coins (Cosmos/Bittensor, TheGraph/InjectiveProtocol) traded places and the #49 row
# (Monero) was replaced by FirstDigitalUSD.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text can never be mis-parsed as a plain number (coin names, URLs,
# "  +x.xx%  " percentages, and multi-dot "thousands.hundreds.cents" prices) -- a
# straight .Value assignment keeps them stored as text, same as the source file.
$textUpdates = @{
    "D2" = '70.626.86'
    "E2" = '  +5.51%  '
    "D3" = '3.617.18'
    "E3" = '  +5.08%  '
    "E4" = '  -0.07%  '
    "E5" = '  +3.59%  '
    "E6" = '  +3.30%  '
    "E7" = '  +1.67%  '
    "D8" = '3.602.07'
    "E8" = '  +4.84%  '
    "E9" = '  -0.09%  '
    "E10" = '  +0.14%  '
    "E12" = '  +5.42%  '
    "E13" = '  +3.76%  '
    "E14" = '  +4.65%  '
    "D15" = '4.184.62'
    "E15" = '  +5.06%  '
    "E16" = '  +4.52%  '
    "D17" = '3.606.57'
    "E17" = '  +4.77%  '
    "D18" = '70.295.69'
    "E18" = '  +5.36%  '
    "E19" = '  +4.06%  '
    "E20" = '  +0.23%  '
    "E21" = '  +4.25%  '
    "E22" = '  +4.02%  '
    "E23" = '  +15.56%  '
    "E24" = '  +7.60%  '
    "E25" = '  +6.73%  '
    "E26" = '  +1.54%  '
    "E27" = '  +5.28%  '
    "E28" = '  +1.76%  '
    "E29" = '  +6.71%  '
    "E30" = '  +2.96%  '
    "E31" = '  +8.74%  '
    "B32" = 'Bittensor'
    "C32" = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
    "E32" = '  +6.46%  '
    "B33" = 'Cosmos'
    "C33" = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
    "E33" = '  +5.46%  '
    "E34" = '  +7.20%  '
    "E35" = '  +3.72%  '
    "D36" = '0.0₃0823'
    "E36" = '  +6.36%  '
    "B37" = 'InjectiveProtocol'
    "C37" = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
    "E37" = '  +4.26%  '
    "B38" = 'TheGraph'
    "C38" = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
    "E38" = '  +3.18%  '
    "E39" = '  +0.09%  '
    "E40" = '  -1.09%  '
    "E41" = '  -0.39%  '
    "D42" = '3.306.44'
    "E43" = '  +6.32%  '
    "E44" = '  +5.14%  '
    "E45" = '  +2.11%  '
    "E46" = '  +2.56%  '
    "E47" = '  +2.05%  '
    "E48" = '  +5.48%  '
    "E49" = '  -2.85%  '
    "E50" = '  +5.31%  '
    "B51" = 'FirstDigitalUSD'
    "C51" = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
    "E51" = '  +0.00%  '
}
foreach ($ref in $textUpdates.Keys) {
    $ws.Range($ref).Value = $textUpdates[$ref]
}

# Price cells whose new text LOOKS like a plain number (e.g. "0.999", "58.36").
# Excel would normally auto-convert those to a numeric value on assignment, but the
# source workbook stores every Price cell as text (t="inlineStr") -- so format each
# as Text first, write the value, then restore the default "Normal" style so no
# stray number formatting is left behind on the cell.
$numericLookingUpdates = @{
    "D4" = '0.999'
    "D5" = '590.61'
    "D6" = '190.91'
    "D7" = '0.645'
    "D9" = '0.999'
    "D10" = '0.178'
    "D11" = '0.662'
    "D12" = '58.36'
    "D13" = '0.0000292'
    "D14" = '9.82'
    "D16" = '19.40'
    "D19" = '12.51'
    "D20" = '0.121'
    "D22" = '494.07'
    "D23" = '17.31'
    "D24" = '5.38'
    "D25" = '4.47'
    "D26" = '90.98'
    "D28" = '11.15'
    "D29" = '9.53'
    "D30" = '32.47'
    "D31" = '7.58'
    "D32" = '626.84'
    "D33" = '12.27'
    "D35" = '65.44'
    "D37" = '38.19'
    "D38" = '0.405'
    "D41" = '3.63'
    "D43" = '3.10'
    "D44" = '0.0446'
    "D45" = '2.69'
    "D46" = '3.31'
    "D48" = '9.15'
    "D49" = '2.72'
    "D50" = '3.30'
    "D51" = '0.999'
}
foreach ($ref in $numericLookingUpdates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $numericLookingUpdates[$ref]
    $cell.Style = "Normal"
}
